$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new columns I (header "I0") and J (header "IF") ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Give the two new header cells the same look (bold/border/centered) as the
# rest of the header row by copying the format from H1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Fill in the I0 / IF values for data rows 2-80 ---
# Row layout: row index (0-based, row 2 = index 0), I value, J value
$ijData = @(
    @(8, 9),
    @(7, 8),
    @(9, 9),
    @(8, 9),
    @(9, 10),
    @(8, 9),
    @(9, 9),
    @(8, 9),
    @(9, 9),
    @(8, 9),
    @(9, 9),
    @(8, 9),
    @(8, 8),
    @(8, 9),
    @(10, 10),
    @(8, 9),
    @(8, 9),
    @(10, 11),
    @(8, 8),
    @(9, 9),
    @(7, 7),
    @(9, 9),
    @(8, 8),
    @(8, 8),
    @(7, 7),
    @(8, 8),
    @(7, 7),
    @(7, 7),
    @(7, 7),
    @(6, 6),
    @(8, 8),
    @(6, 7),
    @(7, 7),
    @(6, 7),
    @(8, 8),
    @(7, 8),
    @(8, 8),
    @(8, 9),
    @(7, 8),
    @(8, 8),
    @(6, 7),
    @(7, 8),
    @(7, 8),
    @(8, 9),
    @(7, 7),
    @(8, 8),
    @(8, 8),
    @(7, 7),
    @(8, 8),
    @(8, 8),
    @(7, 8),
    @(6, 7),
    @(9, 9),
    @(8, 9),
    @(7, 8),
    @(7, 8),
    @(8, 9),
    @(8, 8),
    @(7, 8),
    @(7, 8),
    @(8, 8),
    @(6, 7),
    @(6, 8),
    @(7, 8),
    @(7, 8),
    @(6, 8),
    @(6, 7),
    @(8, 8),
    @(5, 6),
    @(7, 8),
    @(6, 7),
    @(5, 7),
    @(5, 7),
    @(6, 7),
    @(8, 9),
    @(1, 4),
    @(1, 2),
    @(4, 5),
    @(4, 5),
)

$arr = New-Object 'object[,]' 79,2
for ($k = 0; $k -lt $ijData.Count; $k++) {
    $arr[$k, 0] = $ijData[$k][0]
    $arr[$k, 1] = $ijData[$k][1]
}

$ws.Range("I2:J80").Value2 = $arr

Write-Output "applied I0/IF columns"
